$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ViewTemplateList")

# Add the new "Values" header and "Yes,No" data in column E, mirroring the
# formatting already used by column D (Question / Q1).
$ws.Range("D1:D2").Copy()
$ws.Range("E1:E2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E1").Value = "Values"
$ws.Range("E2").Value = "Yes,No"

# Match the new column's width to the rest of the sheet.
$ws.Columns.Item(5).ColumnWidth = 12.7

# Update the active selection to the newly added cell.
$ws.Range("E2").Select()
